$d = $word.ActiveDocument

# wdAlignParagraphJustify = 3
$wdAlignParagraphJustify = 3

$targets = @(
    "comprehensive evaluation",
    "systematically recorded",
    "clear rationale"
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    foreach ($marker in $targets) {
        if ($t -like "*$marker*") {
            $p.Alignment = $wdAlignParagraphJustify
        }
    }
}
